$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C (shifts old C,D,E -> D,E,F)
$ws.Columns("C:C").Insert()

Write-Host "inserted column"

# New "Postal code" column header text
$ws.Range("C10").Value = "Почтовый индекс"

# Row 12 numbers (previously stored as text "2","3","4"; now literal numbers, plus new F12=5)
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 5

